# Update MRIP figures in indicator table
# - D4 (MRIP recreational trips time series image) : 2025-02-18 -> 2025-04-03
# - D5 (MRIP recreational landings time series image) : 2025-02-18 -> 2025-04-03

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = "total_recreational_trips_n_millions_2025-04-03.png"
$ws.Range("D5").Value = "total_recreational_landings_lbs_millions_2025-04-03.png"
